$d = $word.ActiveDocument

# 1. "Available to start for full time employment July 2024" -> "... July 2025"
$d.Content.Find.Execute("July 2024", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "July 2025", 2) | Out-Null

# 2. "Graduating May 2024" -> "Graduating May 2025"
$d.Content.Find.Execute("Graduating May 2024", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Graduating May 2025", 2) | Out-Null

# 3. "Inducted into Tau Beta Pi (November 2022)" -> "Inducted into Tau Beta Pi (April 2024)"
$d.Content.Find.Execute("November 2022", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "April 2024", 2) | Out-Null

# 4. Study Abroad "January - August 2022)" -> "January - August 2023)"
$d.Content.Find.Execute([string]::Concat("January ", [char]0x2013, " August 2022"), $true, $false, $false, $false, $false, `
                         $true, 1, $false, [string]::Concat("January ", [char]0x2013, " August 2023"), 2) | Out-Null

# 5. INITECH Intern "May 2023 - August 2023" - merge the split "202"+"3" run for the first year
#    (text itself is unchanged, this just normalizes formatting/run layout)
$d.Content.Find.Execute("InternMay 2023", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "InternMay 2023", 2) | Out-Null

# 6. THON activity: "Networking Chair (May 2022 - May 2023), Secretary (May 2023 - Present)"
#    -> "Networking Chair (May 2023 - May 2024), Secretary (May 2024 - Present)"
$old6 = [string]::Concat("Networking Chair (May 2022 ", [char]0x2013, " May 2023), Secretary (May 2023 ", [char]0x2013, " Present)")
$new6 = [string]::Concat("Networking Chair (May 2023 ", [char]0x2013, " May 2024), Secretary (May 2024 ", [char]0x2013, " Present)")
$d.Content.Find.Execute($old6, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $new6, 2) | Out-Null

# 7. Penn State Football Team "Coordinator (2022 - Present)" -> "Coordinator (2023 - Present)"
$d.Content.Find.Execute("Coordinator (2022", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Coordinator (2023", 2) | Out-Null
